$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph (by 1-based index) whose text contains the
# given marker string.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($marker) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$marker*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Helper: write $fullText into the (currently empty) paragraph at $paraIndex,
# formatted bold + (bold, complex-script) + dark red (C00000), matching the
# style already used elsewhere in the document for typed-in answers.
# ---------------------------------------------------------------------------
function Set-AnswerText($paraIndex, $fullText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Text = $fullText
    $r.Font.Bold = 1
    $r.Font.BoldBi = 1
    $r.Font.Color = 192
}

# ---------------------------------------------------------------------------
# Helper: force a run split inside the paragraph's text at $paraIndex right
# before $searchText, so the answer ends up encoded as two separate <w:r>
# runs (as produced when text is typed/edited incrementally in Word) while
# keeping identical visible formatting on both pieces.
# ---------------------------------------------------------------------------
function Split-AnswerRun($paraIndex, $searchText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $fr = $p.Range.Duplicate
    $fr.Find.Execute($searchText) | Out-Null
    $fr.Font.Color = 255
    $fr.Font.Color = 192
}

# ---------------------------------------------------------------------------
# Q17 -> 68.27%
# ---------------------------------------------------------------------------
$q17 = Find-ParagraphIndex("Q17. Based on the dataset")
Set-AnswerText ($q17 + 1) "68.27%"

# ---------------------------------------------------------------------------
# Q18 -> 0.9545 (typed as "0.9" then "545")
# ---------------------------------------------------------------------------
$q18 = Find-ParagraphIndex("Q18. Based on the dataset")
Set-AnswerText ($q18 + 1) "0.9545"
Split-AnswerRun ($q18 + 1) "545"

# ---------------------------------------------------------------------------
# Q19 -> Yes
# ---------------------------------------------------------------------------
$q19 = Find-ParagraphIndex("Q19. Are data for length")
Set-AnswerText ($q19 + 1) "Yes"

# ---------------------------------------------------------------------------
# Q20 -> 0.054 (typed as "0.0" then "54")
# ---------------------------------------------------------------------------
$q20 = Find-ParagraphIndex("Q20. What is the significance value")
Set-AnswerText ($q20 + 1) "0.054"
Split-AnswerRun ($q20 + 1) "54"

# ---------------------------------------------------------------------------
# Q21 -> -1.01 (typed as "-" then "1.01")
# ---------------------------------------------------------------------------
$q21 = Find-ParagraphIndex("Q21. What is the standard score")
Set-AnswerText ($q21 + 1) "-1.01"
Split-AnswerRun ($q21 + 1) "1.01"

# ---------------------------------------------------------------------------
# Q22 -> 0.836 (typed as "0." then "836")
# ---------------------------------------------------------------------------
$q22 = Find-ParagraphIndex("Q22. Based on the dataset")
Set-AnswerText ($q22 + 1) "0.836"
Split-AnswerRun ($q22 + 1) "836"

Write-Host "Inserted answers for Q17, Q18, Q19, Q20, Q21, Q22."
